$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Record progress on "List User Ads", "Show Buttons", "Implement Paging",
# "Implemented Category and Filtering" and "Edit Inactive Ads" - each is
# worth the standard "Up to 5" commit score.
$ws.Range("C21").Value = 5
$ws.Range("C22").Value = 5
$ws.Range("C23").Value = 5
$ws.Range("C24").Value = 5
$ws.Range("C25").Value = 5

# Move the view down to where work is happening and leave the selection
# on the newly scored row.
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("G27").Select()
